# Hindalco price sheet update (2025-12-18):
# A new day's price row is published, so every existing data row (2..190)
# shifts down by one (3..191) and a brand new row 2 is inserted at the top
# with the 18-12-2025 figures. Because row 190 (12-06-2025) moves to row
# 191, the table grows from A1:F190 to A1:F191 and a hyperlink needs to be
# (re)created on the newly exposed F130 cell, whose text/URL shifted down
# from the former F129 but whose hyperlink object did not travel with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 2, shifting all data rows down by one.
$ws.Rows("2:2").Insert()

# The freshly inserted row inherited the header row's bold formatting;
# restore the normal data-row formatting by copying it from row 3 (the
# row that used to be row 2, still carrying the correct style).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Populate the new top row with the latest circular's data.
$ws.Range("A2").Value = "18-12-2025"
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 304.5
$ws.Range("E2").Value = "18.12.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-18-december-2025.pdf"

# F130's text (shifted down from the old F129) already shows the
# 12-august-2025 circular URL, but the clickable hyperlink object stayed
# behind on F129. Re-create the hyperlink on F130 pointing at the same
# target used by F129/rId128.
$ws.Hyperlinks.Add($ws.Range("F130"), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf")

# Adding the hyperlink applies default hyperlink formatting; this sheet
# keeps every cell (linked or not) on the plain data style, so restore it
# by copying the format from the neighboring E130 cell.
$ws.Range("E130").Copy()
$ws.Range("F130").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0
